$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 52.55876
$ws.Cells.Item(2, 8).Value = 157.67628
$ws.Cells.Item(2, 9).Value = 0.4767987874074868
$ws.Cells.Item(2, 10).Value = 0.4767987874074869
$ws.Cells.Item(2, 13).Value = 4.043133999999999
$ws.Cells.Item(2, 14).Value = 12.129402
$ws.Cells.Item(2, 15).Value = 0.7517044794313785
$ws.Cells.Item(2, 16).Value = 0.7517044794313784
$ws.Cells.Item(2, 17).Value = 212.50210955384
$ws.Cells.Item(2, 18).Value = 1912.51898598456
$ws.Cells.Item(2, 19).Value = 0.3584117842816574
$ws.Cells.Item(2, 20).Value = 0.3584117842816574
$ws.Cells.Item(3, 7).Value = 52.55876
$ws.Cells.Item(3, 8).Value = 157.67628
$ws.Cells.Item(3, 9).Value = 0.4767987874074868
$ws.Cells.Item(3, 10).Value = 0.4767987874074869
$ws.Cells.Item(3, 15).Value = 0.07962758736516451
$ws.Cells.Item(3, 16).Value = 0.07962758736516451
$ws.Cells.Item(3, 17).Value = 22.51021612453334
$ws.Cells.Item(3, 18).Value = 202.5919451208
$ws.Cells.Item(3, 19).Value = 0.03796633709989415
$ws.Cells.Item(3, 20).Value = 0.03796633709989416
$ws.Cells.Item(4, 7).Value = 52.55876
$ws.Cells.Item(4, 8).Value = 157.67628
$ws.Cells.Item(4, 9).Value = 0.4767987874074868
$ws.Cells.Item(4, 10).Value = 0.4767987874074869
$ws.Cells.Item(4, 13).Value = 0.3522683333333333
$ws.Cells.Item(4, 14).Value = 1.056805
$ws.Cells.Item(4, 15).Value = 0.06549416470700518
$ws.Cells.Item(4, 16).Value = 0.06549416470700517
$ws.Cells.Item(4, 17).Value = 18.51478678726667
$ws.Cells.Item(4, 18).Value = 166.6330810854
$ws.Cells.Item(4, 19).Value = 0.03122753831456629
$ws.Cells.Item(4, 20).Value = 0.03122753831456629
$ws.Cells.Item(5, 7).Value = 52.55876
$ws.Cells.Item(5, 8).Value = 157.67628
$ws.Cells.Item(5, 9).Value = 0.4767987874074868
$ws.Cells.Item(5, 10).Value = 0.4767987874074869
$ws.Cells.Item(5, 13).Value = 0.5549326666666667
$ws.Cells.Item(5, 14).Value = 1.664798
$ws.Cells.Item(5, 15).Value = 0.1031737684964519
$ws.Cells.Item(5, 16).Value = 0.1031737684964519
$ws.Cells.Item(5, 17).Value = 29.16657284349334
$ws.Cells.Item(5, 18).Value = 262.49915559144
$ws.Cells.Item(5, 19).Value = 0.04919312771136902
$ws.Cells.Item(5, 20).Value = 0.04919312771136902
$ws.Cells.Item(6, 9).Value = 0.03798452361347728
$ws.Cells.Item(6, 10).Value = 0.03798452361347729
$ws.Cells.Item(6, 13).Value = 4.043133999999999
$ws.Cells.Item(6, 14).Value = 12.129402
$ws.Cells.Item(6, 15).Value = 0.7517044794313785
$ws.Cells.Item(6, 16).Value = 0.7517044794313784
$ws.Cells.Item(6, 17).Value = 16.929135751688
$ws.Cells.Item(6, 18).Value = 152.362221765192
$ws.Cells.Item(6, 19).Value = 0.02855313654931784
$ws.Cells.Item(6, 20).Value = 0.02855313654931785
$ws.Cells.Item(7, 9).Value = 0.03798452361347728
$ws.Cells.Item(7, 10).Value = 0.03798452361347729
$ws.Cells.Item(7, 15).Value = 0.07962758736516451
$ws.Cells.Item(7, 16).Value = 0.07962758736516451
$ws.Cells.Item(7, 19).Value = 0.003024615972556316
$ws.Cells.Item(7, 20).Value = 0.003024615972556317
$ws.Cells.Item(8, 9).Value = 0.03798452361347728
$ws.Cells.Item(8, 10).Value = 0.03798452361347729
$ws.Cells.Item(8, 13).Value = 0.3522683333333333
$ws.Cells.Item(8, 14).Value = 1.056805
$ws.Cells.Item(8, 15).Value = 0.06549416470700518
$ws.Cells.Item(8, 16).Value = 0.06549416470700517
$ws.Cells.Item(8, 17).Value = 1.474994011086667
$ws.Cells.Item(8, 18).Value = 13.27494609978
$ws.Cells.Item(8, 19).Value = 0.002487764645858208
$ws.Cells.Item(8, 20).Value = 0.002487764645858209
$ws.Cells.Item(9, 9).Value = 0.03798452361347728
$ws.Cells.Item(9, 10).Value = 0.03798452361347729
$ws.Cells.Item(9, 13).Value = 0.5549326666666667
$ws.Cells.Item(9, 14).Value = 1.664798
$ws.Cells.Item(9, 15).Value = 0.1031737684964519
$ws.Cells.Item(9, 16).Value = 0.1031737684964519
$ws.Cells.Item(9, 17).Value = 2.323576326445334
$ws.Cells.Item(9, 18).Value = 20.912186938008
$ws.Cells.Item(9, 19).Value = 0.003919006445744914
$ws.Cells.Item(9, 20).Value = 0.003919006445744915
$ws.Cells.Item(10, 7).Value = 53.437349
$ws.Cells.Item(10, 8).Value = 160.312047
$ws.Cells.Item(10, 9).Value = 0.484769108051078
$ws.Cells.Item(10, 10).Value = 0.4847691080510781
$ws.Cells.Item(10, 13).Value = 4.043133999999999
$ws.Cells.Item(10, 14).Value = 12.129402
$ws.Cells.Item(10, 15).Value = 0.7517044794313785
$ws.Cells.Item(10, 16).Value = 0.7517044794313784
$ws.Cells.Item(10, 17).Value = 216.054362611766
$ws.Cells.Item(10, 18).Value = 1944.489263505894
$ws.Cells.Item(10, 19).Value = 0.3644031100119493
$ws.Cells.Item(10, 20).Value = 0.3644031100119493
$ws.Cells.Item(11, 7).Value = 53.437349
$ws.Cells.Item(11, 8).Value = 160.312047
$ws.Cells.Item(11, 9).Value = 0.484769108051078
$ws.Cells.Item(11, 10).Value = 0.4847691080510781
$ws.Cells.Item(11, 15).Value = 0.07962758736516451
$ws.Cells.Item(11, 16).Value = 0.07962758736516451
$ws.Cells.Item(11, 17).Value = 22.88650407871334
$ws.Cells.Item(11, 18).Value = 205.97853670842
$ws.Cells.Item(11, 19).Value = 0.03860099450327009
$ws.Cells.Item(11, 20).Value = 0.0386009945032701
$ws.Cells.Item(12, 7).Value = 53.437349
$ws.Cells.Item(12, 8).Value = 160.312047
$ws.Cells.Item(12, 9).Value = 0.484769108051078
$ws.Cells.Item(12, 10).Value = 0.4847691080510781
$ws.Cells.Item(12, 13).Value = 0.3522683333333333
$ws.Cells.Item(12, 14).Value = 1.056805
$ws.Cells.Item(12, 15).Value = 0.06549416470700518
$ws.Cells.Item(12, 16).Value = 0.06549416470700517
$ws.Cells.Item(12, 17).Value = 18.82428586998167
$ws.Cells.Item(12, 18).Value = 169.418572829835
$ws.Cells.Item(12, 19).Value = 0.0317495478075653
$ws.Cells.Item(12, 20).Value = 0.0317495478075653
$ws.Cells.Item(13, 7).Value = 53.437349
$ws.Cells.Item(13, 8).Value = 160.312047
$ws.Cells.Item(13, 9).Value = 0.484769108051078
$ws.Cells.Item(13, 10).Value = 0.4847691080510781
$ws.Cells.Item(13, 13).Value = 0.5549326666666667
$ws.Cells.Item(13, 14).Value = 1.664798
$ws.Cells.Item(13, 15).Value = 0.1031737684964519
$ws.Cells.Item(13, 16).Value = 0.1031737684964519
$ws.Cells.Item(13, 17).Value = 29.65413058016734
$ws.Cells.Item(13, 18).Value = 266.887175221506
$ws.Cells.Item(13, 19).Value = 0.05001545572829339
$ws.Cells.Item(13, 20).Value = 0.0500154557282934
$ws.Cells.Item(14, 5).Value = 2
$ws.Cells.Item(14, 6).Value = 0.6666666666666666
$ws.Cells.Item(14, 7).Value = 0.049338
$ws.Cells.Item(14, 8).Value = 0.148014
$ws.Cells.Item(14, 9).Value = 0.0004475809279577863
$ws.Cells.Item(14, 10).Value = 0.0004475809279577865
$ws.Cells.Item(14, 13).Value = 4.043133999999999
$ws.Cells.Item(14, 14).Value = 12.129402
$ws.Cells.Item(14, 15).Value = 0.7517044794313785
$ws.Cells.Item(14, 16).Value = 0.7517044794313784
$ws.Cells.Item(14, 17).Value = 0.199480145292
$ws.Cells.Item(14, 18).Value = 1.795321307628
$ws.Cells.Item(14, 19).Value = 0.0003364485884539211
$ws.Cells.Item(14, 20).Value = 0.0003364485884539212
$ws.Cells.Item(15, 5).Value = 2
$ws.Cells.Item(15, 6).Value = 0.6666666666666666
$ws.Cells.Item(15, 7).Value = 0.049338
$ws.Cells.Item(15, 8).Value = 0.148014
$ws.Cells.Item(15, 9).Value = 0.0004475809279577863
$ws.Cells.Item(15, 10).Value = 0.0004475809279577865
$ws.Cells.Item(15, 15).Value = 0.07962758736516451
$ws.Cells.Item(15, 16).Value = 0.07962758736516451
$ws.Cells.Item(15, 17).Value = 0.02113080756
$ws.Cells.Item(15, 18).Value = 0.19017726804
$ws.Cells.Item(15, 19).Value = 0.00003563978944394004
$ws.Cells.Item(15, 20).Value = 0.00003563978944394004
$ws.Cells.Item(16, 5).Value = 2
$ws.Cells.Item(16, 6).Value = 0.6666666666666666
$ws.Cells.Item(16, 7).Value = 0.049338
$ws.Cells.Item(16, 8).Value = 0.148014
$ws.Cells.Item(16, 9).Value = 0.0004475809279577863
$ws.Cells.Item(16, 10).Value = 0.0004475809279577865
$ws.Cells.Item(16, 13).Value = 0.3522683333333333
$ws.Cells.Item(16, 14).Value = 1.056805
$ws.Cells.Item(16, 15).Value = 0.06549416470700518
$ws.Cells.Item(16, 16).Value = 0.06549416470700517
$ws.Cells.Item(16, 17).Value = 0.01738021503
$ws.Cells.Item(16, 18).Value = 0.15642193527
$ws.Cells.Item(16, 19).Value = 0.00002931393901538148
$ws.Cells.Item(16, 20).Value = 0.00002931393901538149
$ws.Cells.Item(17, 5).Value = 2
$ws.Cells.Item(17, 6).Value = 0.6666666666666666
$ws.Cells.Item(17, 7).Value = 0.049338
$ws.Cells.Item(17, 8).Value = 0.148014
$ws.Cells.Item(17, 9).Value = 0.0004475809279577863
$ws.Cells.Item(17, 10).Value = 0.0004475809279577865
$ws.Cells.Item(17, 13).Value = 0.5549326666666667
$ws.Cells.Item(17, 14).Value = 1.664798
$ws.Cells.Item(17, 15).Value = 0.1031737684964519
$ws.Cells.Item(17, 16).Value = 0.1031737684964519
$ws.Cells.Item(17, 17).Value = 0.027379267908
$ws.Cells.Item(17, 18).Value = 0.2464134111720001
$ws.Cells.Item(17, 19).Value = 0.00004617861104454376
$ws.Cells.Item(17, 20).Value = 0.00004617861104454376
